# Update TPM-derived statistics (G,H,I,J,M,N,O,P,Q,R,S,T) for rows 2-17
# per new TPM recomputation ("update scripts wuth new tpm").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 18.76993966666667
$ws.Range("H2").Value = 56.309819
$ws.Range("I2").Value = 0.1007685501185251
$ws.Range("J2").Value = 0.1007685501185251
$ws.Range("M2").Value = 1.097462
$ws.Range("N2").Value = 3.292386
$ws.Range("O2").Value = 0.3941457306284283
$ws.Range("P2").Value = 0.3941457306284283
$ws.Range("Q2").Value = 20.59929552645933
$ws.Range("R2").Value = 185.393659738134
$ws.Range("S2").Value = 0.03971749381083348
$ws.Range("T2").Value = 0.03971749381083348
$ws.Range("G3").Value = 18.76993966666667
$ws.Range("H3").Value = 56.309819
$ws.Range("I3").Value = 0.1007685501185251
$ws.Range("J3").Value = 0.1007685501185251
$ws.Range("O3").Value = 0.244178053493144
$ws.Range("P3").Value = 0.2441780534931439
$ws.Range("Q3").Value = 12.76151304990967
$ws.Range("R3").Value = 114.853617449187
$ws.Range("S3").Value = 0.02460546842126779
$ws.Range("T3").Value = 0.02460546842126778
$ws.Range("G4").Value = 18.76993966666667
$ws.Range("H4").Value = 56.309819
$ws.Range("I4").Value = 0.1007685501185251
$ws.Range("J4").Value = 0.1007685501185251
$ws.Range("M4").Value = 0.5872626666666667
$ws.Range("N4").Value = 1.761788
$ws.Range("O4").Value = 0.2109112414134909
$ws.Range("P4").Value = 0.2109112414134909
$ws.Range("Q4").Value = 11.02288482181911
$ws.Range("R4").Value = 99.205963396372
$ws.Range("S4").Value = 0.02125322000093571
$ws.Range("T4").Value = 0.0212532200009357
$ws.Range("G5").Value = 18.76993966666667
$ws.Range("H5").Value = 56.309819
$ws.Range("I5").Value = 0.1007685501185251
$ws.Range("J5").Value = 0.1007685501185251
$ws.Range("M5").Value = 0.419791
$ws.Range("N5").Value = 1.259373
$ws.Range("O5").Value = 0.1507649744649369
$ws.Range("P5").Value = 0.1507649744649369
$ws.Range("Q5").Value = 7.879451742609667
$ws.Range("R5").Value = 70.915065683487
$ws.Range("S5").Value = 0.01519236788548815
$ws.Range("T5").Value = 0.01519236788548815
$ws.Range("I6").Value = 0.5130361557055731
$ws.Range("J6").Value = 0.5130361557055731
$ws.Range("M6").Value = 1.097462
$ws.Range("N6").Value = 3.292386
$ws.Range("O6").Value = 0.3941457306284283
$ws.Range("P6").Value = 0.3941457306284283
$ws.Range("Q6").Value = 104.875810703908
$ws.Range("R6").Value = 943.882296335172
$ws.Range("S6").Value = 0.2022110104293733
$ws.Range("T6").Value = 0.2022110104293732
$ws.Range("I7").Value = 0.5130361557055731
$ws.Range("J7").Value = 0.5130361557055731
$ws.Range("O7").Value = 0.244178053493144
$ws.Range("P7").Value = 0.2441780534931439
$ws.Range("S7").Value = 0.1252721698717924
$ws.Range("T7").Value = 0.1252721698717924
$ws.Range("I8").Value = 0.5130361557055731
$ws.Range("J8").Value = 0.5130361557055731
$ws.Range("M8").Value = 0.5872626666666667
$ws.Range("N8").Value = 1.761788
$ws.Range("O8").Value = 0.2109112414134909
$ws.Range("P8").Value = 0.2109112414134909
$ws.Range("Q8").Value = 56.12007364519734
$ws.Range("R8").Value = 505.080662806776
$ws.Range("S8").Value = 0.1082050924898674
$ws.Range("T8").Value = 0.1082050924898674
$ws.Range("I9").Value = 0.5130361557055731
$ws.Range("J9").Value = 0.5130361557055731
$ws.Range("M9").Value = 0.419791
$ws.Range("N9").Value = 1.259373
$ws.Range("O9").Value = 0.1507649744649369
$ws.Range("P9").Value = 0.1507649744649369
$ws.Range("Q9").Value = 40.116123793994
$ws.Range("R9").Value = 361.045114145946
$ws.Range("S9").Value = 0.07734788291454012
$ws.Range("T9").Value = 0.07734788291454012
$ws.Range("G10").Value = 20.061603
$ws.Range("H10").Value = 60.184809
$ws.Range("I10").Value = 0.1077029912330274
$ws.Range("J10").Value = 0.1077029912330274
$ws.Range("M10").Value = 1.097462
$ws.Range("N10").Value = 3.292386
$ws.Range("O10").Value = 0.3941457306284283
$ws.Range("P10").Value = 0.3941457306284283
$ws.Range("Q10").Value = 22.016846951586
$ws.Range("R10").Value = 198.151622564274
$ws.Range("S10").Value = 0.04245067417040881
$ws.Range("T10").Value = 0.04245067417040881
$ws.Range("G11").Value = 20.061603
$ws.Range("H11").Value = 60.184809
$ws.Range("I11").Value = 0.1077029912330274
$ws.Range("J11").Value = 0.1077029912330274
$ws.Range("O11").Value = 0.244178053493144
$ws.Range("P11").Value = 0.2441780534931439
$ws.Range("Q11").Value = 13.639703325273
$ws.Range("R11").Value = 122.757329927457
$ws.Range("S11").Value = 0.02629870675466979
$ws.Range("T11").Value = 0.02629870675466979
$ws.Range("G12").Value = 20.061603
$ws.Range("H12").Value = 60.184809
$ws.Range("I12").Value = 0.1077029912330274
$ws.Range("J12").Value = 0.1077029912330274
$ws.Range("M12").Value = 0.5872626666666667
$ws.Range("N12").Value = 1.761788
$ws.Range("O12").Value = 0.2109112414134909
$ws.Range("P12").Value = 0.2109112414134909
$ws.Range("Q12").Value = 11.781430475388
$ws.Range("R12").Value = 106.032874278492
$ws.Range("S12").Value = 0.02271577158490414
$ws.Range("T12").Value = 0.02271577158490414
$ws.Range("G13").Value = 20.061603
$ws.Range("H13").Value = 60.184809
$ws.Range("I13").Value = 0.1077029912330274
$ws.Range("J13").Value = 0.1077029912330274
$ws.Range("M13").Value = 0.419791
$ws.Range("N13").Value = 1.259373
$ws.Range("O13").Value = 0.1507649744649369
$ws.Range("P13").Value = 0.1507649744649369
$ws.Range("Q13").Value = 8.421680384973001
$ws.Range("R13").Value = 75.79512346475701
$ws.Range("S13").Value = 0.0162378387230447
$ws.Range("T13").Value = 0.0162378387230447
$ws.Range("G14").Value = 51.87415833333333
$ws.Range("H14").Value = 155.622475
$ws.Range("I14").Value = 0.2784923029428744
$ws.Range("J14").Value = 0.2784923029428744
$ws.Range("M14").Value = 1.097462
$ws.Range("N14").Value = 3.292386
$ws.Range("O14").Value = 0.3941457306284283
$ws.Range("P14").Value = 0.3941457306284283
$ws.Range("Q14").Value = 56.92991755281666
$ws.Range("R14").Value = 512.36925797535
$ws.Range("S14").Value = 0.1097665522178128
$ws.Range("T14").Value = 0.1097665522178128
$ws.Range("G15").Value = 51.87415833333333
$ws.Range("H15").Value = 155.622475
$ws.Range("I15").Value = 0.2784923029428744
$ws.Range("J15").Value = 0.2784923029428744
$ws.Range("O15").Value = 0.244178053493144
$ws.Range("P15").Value = 0.2441780534931439
$ws.Range("Q15").Value = 35.26877338340834
$ws.Range("R15").Value = 317.418960450675
$ws.Range("S15").Value = 0.06800170844541403
$ws.Range("T15").Value = 0.06800170844541402
$ws.Range("G16").Value = 51.87415833333333
$ws.Range("H16").Value = 155.622475
$ws.Range("I16").Value = 0.2784923029428744
$ws.Range("J16").Value = 0.2784923029428744
$ws.Range("M16").Value = 0.5872626666666667
$ws.Range("N16").Value = 1.761788
$ws.Range("O16").Value = 0.2109112414134909
$ws.Range("P16").Value = 0.2109112414134909
$ws.Range("Q16").Value = 30.46375655392222
$ws.Range("R16").Value = 274.1738089853
$ws.Range("S16").Value = 0.05873715733778361
$ws.Range("T16").Value = 0.0587371573377836
$ws.Range("G17").Value = 51.87415833333333
$ws.Range("H17").Value = 155.622475
$ws.Range("I17").Value = 0.2784923029428744
$ws.Range("J17").Value = 0.2784923029428744
$ws.Range("M17").Value = 0.419791
$ws.Range("N17").Value = 1.259373
$ws.Range("O17").Value = 0.1507649744649369
$ws.Range("P17").Value = 0.1507649744649369
$ws.Range("Q17").Value = 21.77630480090833
$ws.Range("R17").Value = 195.986743208175
$ws.Range("S17").Value = 0.04198688494186392
$ws.Range("T17").Value = 0.04198688494186392
